$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# E191 was stored as the text "59.0" but should become the numeric value 59
$ws.Range("E191").Value = 59

# New row 193
$ws.Range("A193").Value = "teste5"
$ws.Range("B193").Value = "teste5.com.br"
$ws.Range("C193").Value = "Não padronizado"
$ws.Range("D193").Value = "Sim"
$ws.Range("E193").Value = 59
$ws.Range("F193").Value = "testando"

# New row 194 - E194 must remain text "115.0", not be converted to a number
$ws.Range("A194").Value = "testes65"
$ws.Range("B194").Value = "stees.com.br"
$ws.Range("C194").Value = "Com erro"
$ws.Range("D194").Value = "Não"
$ws.Range("E194").Value = "'115.0"
$ws.Range("F194").Value = "dual"
